# fix bugs in get_dash
# The source data had an extra "Unnamed: 0.1.1.1.1" index column that was
# missing from the exported sheet; also the amount for row 3 and a couple
# of the trailing "staircase" index cells were wrong, and the last source
# row (18) was missing entirely. This script re-creates that shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before the old "F" (type) column ---------------
# This shifts type/amount/date from F/G/H to G/H/I and gives us a fresh,
# blank column F to populate (it inherits style from its left neighbour,
# which already carries the header's bold/border/centered style in row 1).
$ws.Columns("F").Insert()

# --- 2. Header row -----------------------------------------------------------
$ws.Range("F1").Value = "Unnamed: 0.1.1.1.1"

# --- 3. New column F mirrors column E for the data rows that still have a
#        value in E (the "staircase" pattern: each index column runs one
#        row further than the last before turning blank) ------------------
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("F8").Value = 7
$ws.Range("F9").Value = 8
$ws.Range("F10").Value = 9
$ws.Range("F11").Value = 10
$ws.Range("F12").Value = 11
$ws.Range("F13").Value = 12
$ws.Range("F14").Value = 13
$ws.Range("F15").Value = 14
$ws.Range("F16").Value = 15
# F17 stays blank (matches the staircase pattern one row short).

# --- 4. Fix the amount bug in row 4 (was 80, should be 70) ------------------
$ws.Range("H4").Value = 70

# --- 5. Fill in the rest of the staircase that was previously left blank ---
$ws.Range("E17").Value = 16
$ws.Range("C18").Value = 17
$ws.Range("B19").Value = 18

# --- 6. Append the missing trailing record as row 20 ------------------------
$ws.Range("A20").Value = 18
$ws.Range("G20").Value = "clothing"
$ws.Range("H20").Value = 40
$ws.Range("I20").NumberFormat = "@"
$ws.Range("I20").Value = "2021-08-20"

# Give A20 the same header/index style used by the rest of column A, and
# make sure I20 (the new date cell) ends up with the plain/default style
# like every other data cell rather than the "@" text format we used above
# just to stop Excel from auto-converting the date string.
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial($xlPasteFormats)
$ws.Range("H20").Copy()
$ws.Range("I20").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
